# Add a new "Dilution sheet" worksheet with the dilution/sample-weight table,
# and switch the active tab to it.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the header cells on Sheet1 (E1/F1) - text content only changes
# cosmetically (trailing colon -> "Total vol: ") but keep values consistent.
$ws1.Range("E1:F1").Font.Name = "Calibri"
$ws1.Range("E1:F1").Font.Size = 12
$ws1.Range("E1:F1").Font.Color = 0
$ws1.Range("E1").Value = "Sample wt"
$ws1.Range("F1").Value = "Total vol: "

# Add the new "Dilution sheet" worksheet after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Dilution sheet"

$ws2.Range("A1").Value = "Sample"
$ws2.Range("B1").Value = "Vial wt (g)"
$ws2.Range("C1").Value = "Vial wt+Sample (g)"
$ws2.Range("D1").Value = "Sample Wt (g)"
$ws2.Range("E1").Value = "DI added (mL)"
$ws2.Range("F1").Value = "Total vol (mL)"
$ws2.Range("G1").Value = "Vial wt after addition (g)"
$ws2.Range("J1").Value = "*two different glass thickness-> explains differing vial wts"

# Bold, black Calibri 12 across the header row (A1:H1).
$ws2.Range("A1:H1").Font.Name = "Calibri"
$ws2.Range("A1:H1").Font.Size = 12
$ws2.Range("A1:H1").Font.Bold = $true
$ws2.Range("A1:H1").Font.Color = 0

# Regular (non-bold) black Calibri 12 for the trailing note cells.
$ws2.Range("I1:J1").Font.Name = "Calibri"
$ws2.Range("I1:J1").Font.Size = 12
$ws2.Range("I1:J1").Font.Color = 0

# Highlight the primary weight/volume columns in yellow.
$ws2.Range("A1").Interior.Color = 0x99E6FF
$ws2.Range("D1").Interior.Color = 0x99E6FF
$ws2.Range("F1").Interior.Color = 0x99E6FF

$ws2.Range("A1:J1").Select()
$ws2.Activate()
